$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$status = "Ready for handoff"
$overviewDate = "2016-08-20 10:53:16"
$zhcnHandoffDate = "2016-08-20 10:53:12"
$dedeHandoffDate = "2016-08-20 10:53:16"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c97e55e4eefd8060f3f680d18ec1d282603cf76e/e2e/db9b99e9-93b7-4156-a7fd-3b35655629bf.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f1b26c8f21564653c53e1506333cc0efa9eec9c/e2e/db9b99e9-93b7-4156-a7fd-3b35655629bf.md."

# Overview sheet, row 3 (db9b99e9 file)
$overview.Range("E3").Value = $status
$overview.Range("F3").Value = $status
$overview.Range("G3").Value = $overviewDate

# zh-cn sheet, row 3 (db9b99e9 file)
$zhcn.Range("C3").Value = $status
$zhcn.Range("H3").Value = $zhcnHandoffDate
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# de-de sheet, row 3 (db9b99e9 file)
$dede.Range("C3").Value = $status
$dede.Range("H3").Value = $dedeHandoffDate
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
